$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A6").Value = "Discountcurve"
$ws.Range("B6").Value = "USDOIS 31122019"
$ws.Columns.Item(2).ColumnWidth = 15.5
$ws.Range("C10").Select()
